# The only functional change in the target diff is on slide 3: the two
# hyperlinked runs ("現在登入就送彩虹鑰匙" and "!!!") keep their hyperlink
# (rId2 -> SouthSchool.mp4) but lose the extra
# action="ppaction://hlinkfile" attribute on <a:hlinkClick>.
#
# Re-assigning Hyperlink.Address to its current value makes PowerPoint
# regenerate the <a:hlinkClick> element from the Hyperlink object alone,
# which naturally omits the (redundant, PowerPoint-generated) "open file"
# action URI while preserving the r:id relationship.
#
# Walk every shape/run in the deck so the fix is general rather than
# hard-coded to a slide/shape index.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $runCount = $tr.Runs().Count
                for ($ri = 1; $ri -le $runCount; $ri++) {
                    $run = $tr.Runs($ri)
                    $hl = $run.ActionSettings(1).Hyperlink
                    if ($hl.Address -ne "") {
                        $hl.Address = $hl.Address
                    }
                }
            }
        }
    }
}
